# "update salt prices and loadings in all scenarios"
#
# Sheet1 holds a table of TEA distribution parameters. Row 8 is the
# "Magnesium chloride unit price" scenario and row 9 is the
# "Zinc sulfate unit price" scenario. For both rows:
#   - column E (Lower) gets a new loaded value
#   - columns G (Midpoint) and I (Upper) were previously computed via
#     formulas (=E*0.08 / =E*1.2 and similar) but are now overwritten
#     with plain literal numbers (the values the user pasted in).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: Magnesium chloride unit price
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9: Zinc sulfate unit price
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# The author's last on-sheet selection ended up on rows 8:9 (A8, block A8:XFD9)
$ws.Range("A8:XFD9").Select()
